# Pavani Gandepalli status on 28/1/2021
# Fills in the status-tracker rows for 22/1, 25/1 (on leave), 27/1 and 28/1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 18 - 22/1/2021
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "22/1/2021"
$ws.Range("A18").HorizontalAlignment = -4131

$ws.Range("B18").Value = "1. I have learnt syntax, if,if-else,elif ladder, for loop, basics in list and tuple as my interview profile needed python `n2. Revision done on GIT and ADB commands `n3. Attended the interview at 5PM and shared the interview question to the team"
$ws.Range("B18").WrapText = $true

$ws.Rows.Item(18).RowHeight = 90

# ---------------------------------------------------------------------------
# Row 19 - On leave 25/1/2021 (styled with the "Neutral" cell style, left
# aligned, matching the highlighted "on leave" rows used elsewhere in the
# sheet e.g. A7).
# ---------------------------------------------------------------------------
$ws.Range("A7").Copy($ws.Range("A19"))
$ws.Range("A19").WrapText = $false
$ws.Range("A19").HorizontalAlignment = -4131
$ws.Range("A19").Value = "On leave 25/1/2021"

# ---------------------------------------------------------------------------
# Row 20 - 27/1/2021
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = "27/1/2021"
$ws.Range("A20").HorizontalAlignment = -4131

$ws.Range("D20").Value = "staircase.txt"

$ws.Range("B20").Value = "1. Practiced different storage classes, and checked the size using size<executable_file> to check where the variable is storing `n2. Worked on hacker rank program it took more time to get logic `n3. Understand Bitwise operators and Completed 3 programs in Bitwise operators from givem list of programs "
$ws.Range("B20").WrapText = $true

$ws.Rows.Item(20).RowHeight = 120

# ---------------------------------------------------------------------------
# Row 21 - 28/1/2021
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "28/1/2021"
$ws.Range("A21").HorizontalAlignment = -4131

$ws.Range("B21").Value = "1. Completed 5 programs in bitwise operators from given list`n2. Attended testing session at 10AM and Srinivas gave few tasks on testcases `n3. Completed the task given by Srinivas and attended one more session at 2PM`n4. Completed on Hacker rank program in C"
$ws.Range("B21").WrapText = $true

$ws.Range("D21").Value = "MinMax.txt"

$ws.Rows.Item(21).RowHeight = 105

# ---------------------------------------------------------------------------
# Update the view: scrolled down a couple more rows, selection now on E21.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("E21").Select()
